$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8 through 30 (old extra rows), keep row 1 header and rows 2-7 which get rewritten below.
$ws.Range("A8:A30").EntireRow.Delete() | Out-Null

$ws.Range("A2").Value = "('Angel', ['Token Creature — Angel', 'Flying', '4/4'])"
$ws.Range("A3").Value = "('Demon', ['Token Creature — Demon', 'Flying', '5/5'])"
$ws.Range("A4").Value = "('Human', ['Token Creature — Human', '1/1'])"
$ws.Range("A5").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A6").Value = "('Tamiyo, the Moon Sage Emblem', ['Emblem — Tamiyo', 'You have no maximum hand size.', 'Whenever a card is put into your graveyard from anywhere, you may return it to your hand.'])"
$ws.Range("A7").Value = "('Zombie', ['Token Creature — Zombie', '2/2'])"
